$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$ws.Cells.Item(2, 4).Value = '26.328.54'
$ws.Cells.Item(2, 5).Value = '  -1.97%  '
$ws.Cells.Item(3, 4).Value = '1.835.93'
$ws.Cells.Item(3, 5).Value = '  -2.23%  '
Set-TextValue $ws.Cells.Item(4, 4) '1.001'
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
Set-TextValue $ws.Cells.Item(5, 4) '258.20'
$ws.Cells.Item(5, 5).Value = '  -7.36%  '
Set-TextValue $ws.Cells.Item(6, 4) '1.001'
$ws.Cells.Item(6, 5).Value = '  +0.06%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.5200'
$ws.Cells.Item(7, 5).Value = '  -1.20%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.3225'
$ws.Cells.Item(8, 5).Value = '  -6.41%  '
Set-TextValue $ws.Cells.Item(9, 4) '0.06741'
$ws.Cells.Item(9, 5).Value = '  -2.96%  '
Set-TextValue $ws.Cells.Item(10, 4) '18.45'
$ws.Cells.Item(10, 5).Value = '  -8.36%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.7588'
$ws.Cells.Item(11, 5).Value = '  -5.86%  '
Set-TextValue $ws.Cells.Item(12, 4) '0.07642'
$ws.Cells.Item(12, 5).Value = '  -2.71%  '
$ws.Cells.Item(13, 4).Value = '1.831.73'
$ws.Cells.Item(13, 5).Value = '  -1.13%  '
Set-TextValue $ws.Cells.Item(14, 4) '88.30'
$ws.Cells.Item(14, 5).Value = '  -1.83%  '
Set-TextValue $ws.Cells.Item(15, 4) '5.010'
$ws.Cells.Item(15, 5).Value = '  -2.91%  '
Set-TextValue $ws.Cells.Item(16, 4) '1.001'
$ws.Cells.Item(16, 5).Value = '  +0.08%  '
Set-TextValue $ws.Cells.Item(17, 4) '13.90'
$ws.Cells.Item(17, 5).Value = '  -4.66%  '
Set-TextValue $ws.Cells.Item(18, 4) '1.001'
$ws.Cells.Item(18, 5).Value = '  -0.01%  '
Set-TextValue $ws.Cells.Item(19, 4) '0.000007873'
$ws.Cells.Item(19, 5).Value = '  -2.40%  '
$ws.Cells.Item(20, 4).Value = '26.361.85'
$ws.Cells.Item(20, 5).Value = '  -2.01%  '
$ws.Cells.Item(21, 4).Value = '2.071.01'
$ws.Cells.Item(21, 5).Value = '  -3.69%  '
Set-TextValue $ws.Cells.Item(22, 4) '4.538'
$ws.Cells.Item(22, 5).Value = '  -4.51%  '
Set-TextValue $ws.Cells.Item(23, 4) '9.406'
$ws.Cells.Item(23, 5).Value = '  -6.13%  '
Set-TextValue $ws.Cells.Item(24, 4) '5.912'
Set-TextValue $ws.Cells.Item(25, 4) '144.11'
$ws.Cells.Item(25, 5).Value = '  -1.58%  '
Set-TextValue $ws.Cells.Item(26, 4) '2.229'
$ws.Cells.Item(26, 5).Value = '  -4.86%  '
Set-TextValue $ws.Cells.Item(27, 4) '1.644'
$ws.Cells.Item(27, 5).Value = '  -1.47%  '
$ws.Cells.Item(28, 5).Value = '  -2.72%  '
Set-TextValue $ws.Cells.Item(29, 4) '111.32'
$ws.Cells.Item(29, 5).Value = '  -2.12%  '
Set-TextValue $ws.Cells.Item(30, 4) '4.162'
Set-TextValue $ws.Cells.Item(31, 4) '4.126'
$ws.Cells.Item(31, 5).Value = '  -4.64%  '
Set-TextValue $ws.Cells.Item(32, 4) '0.08705'
$ws.Cells.Item(32, 5).Value = '  -2.37%  '
Set-TextValue $ws.Cells.Item(33, 4) '0.04766'
$ws.Cells.Item(33, 5).Value = '  -3.49%  '
Set-TextValue $ws.Cells.Item(34, 4) '2.854'
$ws.Cells.Item(34, 5).Value = '  -1.36%  '
$ws.Cells.Item(35, 5).Value = '  -5.41%  '
Set-TextValue $ws.Cells.Item(36, 4) '0.6949'
$ws.Cells.Item(36, 5).Value = '  -5.47%  '
Set-TextValue $ws.Cells.Item(37, 4) '3.058'
$ws.Cells.Item(37, 5).Value = '  -6.70%  '
Set-TextValue $ws.Cells.Item(38, 4) '0.01754'
$ws.Cells.Item(38, 5).Value = '  -5.20%  '
Set-TextValue $ws.Cells.Item(39, 4) '2.196'
$ws.Cells.Item(39, 5).Value = '  -8.56%  '
Set-TextValue $ws.Cells.Item(40, 4) '0.4822'
$ws.Cells.Item(40, 5).Value = '  -6.13%  '
Set-TextValue $ws.Cells.Item(41, 4) '110.84'
$ws.Cells.Item(41, 5).Value = '  -4.90%  '
$ws.Cells.Item(42, 5).Value = '  -1.96%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.8817'
$ws.Cells.Item(43, 5).Value = '  -8.21%  '
Set-TextValue $ws.Cells.Item(44, 4) '1.001'
$ws.Cells.Item(44, 5).Value = '  +0.11%  '
Set-TextValue $ws.Cells.Item(45, 4) '7.642'
$ws.Cells.Item(45, 5).Value = '  -5.41%  '
Set-TextValue $ws.Cells.Item(46, 4) '0.4121'
$ws.Cells.Item(46, 5).Value = '  -8.50%  '
Set-TextValue $ws.Cells.Item(47, 4) '0.05842'
$ws.Cells.Item(47, 5).Value = '  -1.72%  '
Set-TextValue $ws.Cells.Item(48, 4) '8.982'
$ws.Cells.Item(48, 5).Value = '  -4.17%  '
Set-TextValue $ws.Cells.Item(49, 4) '0.1228'
$ws.Cells.Item(49, 5).Value = '  -8.58%  '
Set-TextValue $ws.Cells.Item(50, 4) '34.58'
$ws.Cells.Item(50, 5).Value = '  -5.50%  '
Set-TextValue $ws.Cells.Item(51, 4) '0.8812'
$ws.Cells.Item(51, 5).Value = '  -0.30%  '
